# Update temperature (average_county_temperature, column I) and the
# dependent worst/best ASHP COP values (columns N and O) with refreshed
# NOAA-derived figures for the affected facility rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 4 / 5 (facility_id 1000183, reporting rows 2 & 3)
$ws.Range("I4").Value = 21.28240740740739
$ws.Range("N4").Value = 1.368667816644515
$ws.Range("O4").Value = 1.441442038370027
$ws.Range("I5").Value = 21.28240740740739

# Row 14 / 15 (facility_id 1000961)
$ws.Range("I14").Value = 3.38888888888889
$ws.Range("N14").Value = 1.233817681248088
$ws.Range("O14").Value = 1.291146001942376
$ws.Range("I15").Value = 3.38888888888889

# Row 26 / 27 (facility_id 1003048)
$ws.Range("I26").Value = -1.226851851851833
$ws.Range("N26").Value = 1.203236793039155
$ws.Range("O26").Value = 1.257328254301852
$ws.Range("I27").Value = -1.226851851851833

# Row 35 / 36 (facility_id 1004199)
$ws.Range("I35").Value = 21.79166666666666
$ws.Range("N35").Value = 1.372938473321419
$ws.Range("O35").Value = 1.446233342398694
$ws.Range("I36").Value = 21.79166666666666

# Row 37 / 38 (facility_id 1004396)
$ws.Range("I37").Value = 21.19907407407406
$ws.Range("N37").Value = 1.367971510132557
$ws.Range("O37").Value = 1.440661027663225
$ws.Range("I38").Value = 21.19907407407406

# Row 39 / 40 (facility_id 1004874)
$ws.Range("I39").Value = 21.19907407407406
$ws.Range("N39").Value = 1.367971510132557
$ws.Range("O39").Value = 1.440661027663225
$ws.Range("I40").Value = 21.19907407407406

# Row 49 / 50 (facility_id 1006366)
$ws.Range("I49").Value = -1.226851851851833
$ws.Range("N49").Value = 1.203236793039155
$ws.Range("O49").Value = 1.257328254301852
$ws.Range("I50").Value = -1.226851851851833

# Row 57 / 58 / 59 (facility_id 1007135)
$ws.Range("I57").Value = 13.17361111111111
$ws.Range("I58").Value = 13.17361111111111
$ws.Range("N58").Value = 1.304077921028169
$ws.Range("O58").Value = 1.369214264257821
$ws.Range("I59").Value = 13.17361111111111
